$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename D1, add new E1 ---
$ws.Range("D1").Value = "Processing"
$ws.Range("E1").Value = "Recipes"
# Match the header styling (bold, bordered, centered/top) applied to the other header cells
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").VerticalAlignment = -4160
$ws.Range("E1").Borders.LineStyle = 1

# --- Row 2: strip the "Stat/EffectAmount" prefix from B2, add empty E2 ---
$ws.Range("B2").Value = "ATK %`n11`nStability %`n10`nPhysical Pierce %`n20`nASPD`n900`n%`nstronger against Light`n10`nDark Element`n0`nGuard Break %30"
$ws.Range("E2").WrapText = $false

# --- Row 3 (new): duplicate Anguish Sword entry, NPC source ---
$ws.Range("A3").Value = "Anguish Sword[1 Handed Sword]"
$ws.Range("B3").Value = "ATK %`n11`nStability %`n10`nPhysical Pierce %`n20`nASPD`n900`n%`nstronger against Light`n10`nDark Element`n0`nGuard Break %30"
$ws.Range("C3").Value = "[NPC]Blacksmith : ZaldoSofya City: Blacksmith"
$ws.Range("D3").Value = "Sell0 SpinaProcessunknown"
$ws.Range("E3").WrapText = $false

# --- Row 4 (new): duplicate Anguish Sword entry, empty source ---
$ws.Range("A4").Value = "Anguish Sword[1 Handed Sword]"
$ws.Range("B4").Value = "ATK %`n11`nStability %`n10`nPhysical Pierce %`n20`nASPD`n900`n%`nstronger against Light`n10`nDark Element`n0`nGuard Break %30"
$ws.Range("C4").Value = "empty"
$ws.Range("D4").Value = "Sell0 SpinaProcessunknown"
$ws.Range("E4").WrapText = $false

# --- Row 5 (new): Accordion-fold Sword ---
$ws.Range("A5").Value = "Accordion-fold Sword[1 Handed Sword]"
$ws.Range("B5").Value = "Base ATK`n1`nAggro %`n30`nAttack MP Recovery`n3`nBase Stability %40"
$ws.Range("C5").Value = "[NPC]Blacksmith : ZaldoSofya City: Blacksmith"
$ws.Range("D5").Value = "SellUnknownProcessN/A"
$ws.Range("E5").Value = "Fee`n150`n SpinaSet`n1`n pcsLevelN/ADifficulty`n0`nMaterials- `n25`nxBat Wing- `n1`nxFour-leaf Clover- `n50`nx Cloth- `n25`nx Metal"

# --- Row 6 (new): Shortsword, Pico monster drop ---
$ws.Range("A6").Value = "Shortsword[1 Handed Sword]"
$ws.Range("B6").Value = "Base ATK`n10`nBase Stability %80"
$ws.Range("C6").Value = "Pico(Lv 1)Rakau Plains"
$ws.Range("D6").Value = "Sell1 SpinaProcess2 Metal"
$ws.Range("E6").WrapText = $false

# --- Row 7 (new): Shortsword, NPC source ---
$ws.Range("A7").Value = "Shortsword[1 Handed Sword]"
$ws.Range("B7").Value = "Base ATK`n10`nBase Stability %80"
$ws.Range("C7").Value = "[NPC]Blacksmith : ZaldoSofya City: Blacksmith"
$ws.Range("D7").Value = "Sell1 SpinaProcess2 Metal"
$ws.Range("E7").Value = "Fee`n5`n SpinaSet`n1`n pcsLevel`n1`nDifficulty`n1`nMaterials- `n3`nx Metal"

# --- Row 8 (new): Shortsword, Player Smith Skill ---
$ws.Range("A8").Value = "Shortsword[1 Handed Sword]"
$ws.Range("B8").Value = "Base ATK`n10`nBase Stability %`n80`nPotential15"
$ws.Range("C8").Value = "[Player]Smith Skill-"
$ws.Range("D8").Value = "Sell1 SpinaProcess2 Metal"
$ws.Range("E8").Value = "Fee`n5`n SpinaSet`n1`n pcsLevel`n1`nDifficulty`n1`nMaterials- `n3`nx Metal"

# --- Row 9 (new): Wood Sword ---
$ws.Range("A9").Value = "Wood Sword[1 Handed Sword]"
$ws.Range("B9").Value = "Base ATK`n10`nBase Stability %`n40`nPotential15"
$ws.Range("C9").Value = "[Player]Smith Skill-"
$ws.Range("D9").Value = "Sell1 SpinaProcess2 Wood"
$ws.Range("E9").Value = "FeeN/ASet`n1`n pcsLevel`n1`nDifficulty`n0`nMaterials- `n25`nx Wood"

# --- Row 10 (new): Longsword, NPC source ---
$ws.Range("A10").Value = "Longsword[1 Handed Sword]"
$ws.Range("B10").Value = "Base ATK`n17`nMaxHP`n50`nAccuracy`n1`nBase Stability %80"
$ws.Range("C10").Value = "[NPC]Blacksmith : ZaldoSofya City: Blacksmith"
$ws.Range("D10").Value = "Sell10 SpinaProcess5 Metal"
$ws.Range("E10").Value = "Fee`n50`n SpinaSet`n1`n pcsLevel`n5`nDifficulty`n10`nMaterials- `n20`nxSmall Hilt- `n25`nx Metal"

# --- Row 11 (new): Longsword, Player Smith Skill ---
$ws.Range("A11").Value = "Longsword[1 Handed Sword]"
$ws.Range("B11").Value = "Base ATK`n17`nBase Stability %`n80`nPotential16"
$ws.Range("C11").Value = "[Player]Smith Skill-"
$ws.Range("D11").Value = "Sell10 SpinaProcess5 Metal"
$ws.Range("E11").Value = "Fee`n50`n SpinaSet`n1`n pcsLevel`n5`nDifficulty`n10`nMaterials- `n20`nxSmall Hilt- `n25`nx Metal"

# --- Row 12 (new): Gladius ---
$ws.Range("A12").Value = "Gladius[1 Handed Sword]"
$ws.Range("B12").Value = "Base ATK`n25`nASPD %`n5`nCritical Rate`n1`nBase Stability %80"
$ws.Range("C12").Value = "Goblin(Lv 20)Ruined Temple: Area 1Goblin(Lv 20)Ruined Temple: Area 2"
$ws.Range("D12").Value = "Sell20 SpinaProcess10 Metal"
$ws.Range("E12").WrapText = $false
